# Applies the commit's spreadsheet edits:
#  - DataBase sheet: move selection from B6 to B2
#  - SQLite sheet: move selection from I11 to G20
#  - SQLite sheet: tighten up type labels (String -> Boolean, int -> Integer)

$wb = $excel.ActiveWorkbook

$wsDataBase = $wb.Worksheets.Item("DataBase")
$wsDataBase.Activate()
$wsDataBase.Range("B2").Select()

$wsSQLite = $wb.Worksheets.Item("SQLite")
$wsSQLite.Activate()
$wsSQLite.Range("C3").Value = "Boolean"
$wsSQLite.Range("E3").Value = "Integer"
$wsSQLite.Range("G20").Select()
